$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").ClearContents() | Out-Null

# Row 3
$ws.Range("K3").Value2 = -2.184539318894706
$ws.Range("J3").Value2 = -2.419547318894706
$ws.Range("I3").Value2 = -1.489825318894706
$ws.Range("H3").Value2 = -4.004270318894706
$ws.Range("G3").Value2 = -3.708941318894706
$ws.Range("F3").Value2 = -0.5308223188947059
$ws.Range("E3").Value2 = -2.319131318894706
$ws.Range("D3").Value2 = -10.45921331889471
$ws.Range("C3").Value2 = 7.869792681105293
$ws.Range("B3").Value2 = [double]"1.052939957446597E-10"

# Row 4
$ws.Range("K4").Value2 = -9.960825649752
$ws.Range("J4").Value2 = -9.653689649752
$ws.Range("I4").Value2 = -9.888697649751998
$ws.Range("H4").Value2 = -8.958975649751999
$ws.Range("G4").Value2 = -11.473420649752
$ws.Range("F4").Value2 = -11.178091649752
$ws.Range("E4").Value2 = -7.999972649751999
$ws.Range("D4").Value2 = -9.788281649751999
$ws.Range("C4").Value2 = -17.928363649752
$ws.Range("B4").Value2 = 0.4006423502480008

# Row 5
$ws.Range("K5").Value2 = 8.721112150385469
$ws.Range("J5").Value2 = 8.210436150385469
$ws.Range("I5").Value2 = 8.517572150385469
$ws.Range("H5").Value2 = 8.28256415038547
$ws.Range("G5").Value2 = 9.21228615038547
$ws.Range("F5").Value2 = 6.69784115038547
$ws.Range("E5").Value2 = 6.99317015038547
$ws.Range("D5").Value2 = 10.17128915038547
$ws.Range("C5").Value2 = 8.382980150385469
$ws.Range("B5").Value2 = 0.2428981503854697

# Row 6
$ws.Range("K6").Value2 = 0.2940736034063205
$ws.Range("J6").Value2 = 0.5814406034063205
$ws.Range("I6").Value2 = 0.07076460340632051
$ws.Range("H6").Value2 = 0.3779006034063205
$ws.Range("G6").Value2 = 0.1428926034063205
$ws.Range("F6").Value2 = 1.072614603406321
$ws.Range("E6").Value2 = -1.441830396593679
$ws.Range("D6").Value2 = -1.146501396593679
$ws.Range("C6").Value2 = 2.03161760340632
$ws.Range("B6").Value2 = 0.2433086034063205

# Row 7
$ws.Range("K7").Value2 = -1.825497868393072
$ws.Range("J7").Value2 = -1.853496868393072
$ws.Range("I7").Value2 = -1.566129868393072
$ws.Range("H7").Value2 = -2.076805868393072
$ws.Range("G7").Value2 = -1.769669868393072
$ws.Range("F7").Value2 = -2.004677868393072
$ws.Range("E7").Value2 = -1.074955868393072
$ws.Range("D7").Value2 = -3.589400868393072
$ws.Range("C7").Value2 = -3.294071868393072
$ws.Range("B7").Value2 = -0.115952868393072

# Row 8
$ws.Range("K8").Value2 = 1.707180118130258
$ws.Range("J8").Value2 = 1.614532118130258
$ws.Range("I8").Value2 = 1.586533118130258
$ws.Range("H8").Value2 = 1.873900118130258
$ws.Range("G8").Value2 = 1.363224118130258
$ws.Range("F8").Value2 = 1.670360118130258
$ws.Range("E8").Value2 = 1.435352118130258
$ws.Range("D8").Value2 = 2.365074118130258
$ws.Range("C8").Value2 = -0.1493708818697419
$ws.Range("B8").Value2 = 0.1459581181302581

# Row 9
$ws.Range("K9").Value2 = 1.981330296544459
$ws.Range("J9").Value2 = 1.768354296544459
$ws.Range("I9").Value2 = 1.675706296544459
$ws.Range("H9").Value2 = 1.647707296544459
$ws.Range("G9").Value2 = 1.935074296544459
$ws.Range("F9").Value2 = 1.424398296544459
$ws.Range("E9").Value2 = 1.731534296544459
$ws.Range("D9").Value2 = 1.496526296544459
$ws.Range("C9").Value2 = 2.426248296544459
$ws.Range("B9").Value2 = -0.08819670345554087

# Row 10
$ws.Range("K10").Value2 = -0.5607856406117955
$ws.Range("J10").Value2 = -0.06012564061179543
$ws.Range("I10").Value2 = -0.2731016406117954
$ws.Range("H10").Value2 = -0.3657496406117954
$ws.Range("G10").Value2 = -0.3937486406117954
$ws.Range("F10").Value2 = -0.1063816406117954
$ws.Range("E10").Value2 = -0.6170576406117954
$ws.Range("D10").Value2 = -0.3099216406117954
$ws.Range("C10").Value2 = -0.5449296406117954
$ws.Range("B10").Value2 = 0.3847923593882046

# Row 11
$ws.Range("K11").Value2 = 0.3115658964218456
$ws.Range("J11").Value2 = 0.02180889642184558
$ws.Range("I11").Value2 = 0.5224688964218456
$ws.Range("H11").Value2 = 0.3094928964218456
$ws.Range("G11").Value2 = 0.2168448964218456
$ws.Range("F11").Value2 = 0.1888458964218456
$ws.Range("E11").Value2 = 0.4762128964218456
$ws.Range("D11").Value2 = -0.03446310357815441
$ws.Range("C11").Value2 = 0.2726728964218456
$ws.Range("B11").Value2 = 0.03766489642184559

# Row 12
$ws.Range("K12").Value2 = -0.04200114438751221
$ws.Range("J12").Value2 = 0.2158908556124878
$ws.Range("I12").Value2 = -0.07386614438751221
$ws.Range("H12").Value2 = 0.4267938556124878
$ws.Range("G12").Value2 = 0.2138178556124878
$ws.Range("F12").Value2 = 0.1211698556124878
$ws.Range("E12").Value2 = 0.09317085561248779
$ws.Range("D12").Value2 = 0.3805378556124878
$ws.Range("C12").Value2 = -0.1301381443875122
$ws.Range("B12").Value2 = 0.1769978556124878

# Row 13
$ws.Range("K13").Value2 = 0.6007988926112107
$ws.Range("J13").Value2 = 0.2184178926112106
$ws.Range("I13").Value2 = 0.4763098926112106
$ws.Range("H13").Value2 = 0.1865528926112106
$ws.Range("G13").Value2 = 0.6872128926112107
$ws.Range("F13").Value2 = 0.4742368926112106
$ws.Range("E13").Value2 = 0.3815888926112106
$ws.Range("D13").Value2 = 0.3535898926112106
$ws.Range("C13").Value2 = 0.6409568926112106
$ws.Range("B13").Value2 = 0.1302808926112106

# Row 14
$ws.Range("K14").Value2 = -0.6459331975472806
$ws.Range("J14").Value2 = -0.2346561035472806
$ws.Range("I14").Value2 = -0.6170371035472806
$ws.Range("H14").Value2 = -0.3591451035472806
$ws.Range("G14").Value2 = -0.6489021035472806
$ws.Range("F14").Value2 = -0.1482421035472806
$ws.Range("E14").Value2 = -0.3612181035472806
$ws.Range("D14").Value2 = -0.4538661035472806
$ws.Range("C14").Value2 = -0.4818651035472806
$ws.Range("B14").Value2 = -0.1944981035472806

# Row 15
$ws.Range("K15").Value2 = -0.2588977436446591
$ws.Range("J15").Value2 = -0.5458538376446591
$ws.Range("I15").Value2 = -0.1345767436446591
$ws.Range("H15").Value2 = -0.5169577436446591
$ws.Range("G15").Value2 = -0.2590657436446591
$ws.Range("F15").Value2 = -0.5488227436446591
$ws.Range("E15").Value2 = -0.04816274364465911
$ws.Range("D15").Value2 = -0.2611387436446591
$ws.Range("C15").Value2 = -0.3537867436446591
$ws.Range("B15").Value2 = -0.3817857436446591

# Row 16
$ws.Range("J16").Value2 = 0.09488876243503713
$ws.Range("I16").Value2 = -0.1920673315649629
$ws.Range("H16").Value2 = 0.2192097624350371
$ws.Range("G16").Value2 = -0.1631712375649629
$ws.Range("F16").Value2 = 0.09472076243503715
$ws.Range("E16").Value2 = -0.1950362375649629
$ws.Range("D16").Value2 = 0.3056237624350371
$ws.Range("C16").Value2 = 0.09264776243503714
$ws.Range("B16").Value2 = [double]"-2.375649628613696E-07"

# Row 17
$ws.Range("I17").Value2 = -0.0532886279974082
$ws.Range("H17").Value2 = -0.3402447219974082
$ws.Range("G17").Value2 = 0.0710323720025918
$ws.Range("F17").Value2 = -0.3113486279974082
$ws.Range("E17").Value2 = -0.05345662799740819
$ws.Range("D17").Value2 = -0.3432136279974082
$ws.Range("C17").Value2 = 0.1574463720025918
$ws.Range("B17").Value2 = -0.0555296279974082

# Row 18
$ws.Range("H18").Value2 = -0.2107346170015632
$ws.Range("G18").Value2 = -0.4976907110015631
$ws.Range("F18").Value2 = -0.08641361700156319
$ws.Range("E18").Value2 = -0.4687946170015632
$ws.Range("D18").Value2 = -0.2109026170015632
$ws.Range("C18").Value2 = -0.5006596170015631
$ws.Range("B18").Value2 = [double]"3.829984367986761E-07"

# Row 19
$ws.Range("G19").Value2 = 0.2899248395245076
$ws.Range("F19").Value2 = 0.002968745524507627
$ws.Range("E19").Value2 = 0.4142458395245076
$ws.Range("D19").Value2 = 0.03186483952450761
$ws.Range("C19").Value2 = 0.2897568395245076
$ws.Range("B19").Value2 = [double]"-1.604754923945073E-07"

# Row 20
$ws.Range("F20").Value2 = 0.009566958989038449
$ws.Range("E20").Value2 = -0.2773891350109615
$ws.Range("D20").Value2 = 0.1338879589890384
$ws.Range("C20").Value2 = -0.2484930410109615
$ws.Range("B20").Value2 = 0.009398958989038461

# Row 21
$ws.Range("E21").Value2 = 0.1815418168340501
$ws.Range("D21").Value2 = -0.1054142771659499
$ws.Range("C21").Value2 = 0.3058628168340501
$ws.Range("B21").Value2 = -0.07651818316594991

# Row 22
$ws.Range("D22").Value2 = -0.1243207229122814
$ws.Range("C22").Value2 = -0.4112768169122814
$ws.Range("B22").Value2 = [double]"2.770877186031306E-07"

# Row 23
$ws.Range("C23").Value2 = 0.4880092297750048
$ws.Range("B23").Value2 = 0.2010531357750048

# Row 24
$ws.Range("B24").Value2 = -0.2003621554241067
